# Weekly update: insert two new daily price records for "Perejil" right
# after the existing row 100, shifting the remaining historical rows
# down by two (old row 101 -> 103, ..., old row 123 -> 125).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two new records.
$ws.Rows("101:102").Insert()

# New record 1 (Primera quality, market date 2023-10-05).
$ws.Range("A101").Value = 7
$ws.Range("B101").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C101").Value = "Ñuble"
$ws.Range("D101").Value = 45204
$ws.Range("E101").Value = 16
$ws.Range("F101").Value = 100112044
$ws.Range("G101").Value = "Perejil"
$ws.Range("H101").Value = "Sin especificar"
$ws.Range("I101").Value = "Primera"
$ws.Range("J101").Value = 300
$ws.Range("K101").Value = 2000
$ws.Range("L101").Value = 2000
$ws.Range("M101").Value = 2000
$ws.Range("N101").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O101").Value = "Región de Ñuble"
$ws.Range("P101").Value = 2000
$ws.Range("Q101").Value = 1
$ws.Range("R101").Value = "Hortaliza"

# New record 2 (Segunda quality, market date 2023-10-05).
$ws.Range("A102").Value = 7
$ws.Range("B102").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C102").Value = "Ñuble"
$ws.Range("D102").Value = 45204
$ws.Range("E102").Value = 16
$ws.Range("F102").Value = 100112044
$ws.Range("G102").Value = "Perejil"
$ws.Range("H102").Value = "Sin especificar"
$ws.Range("I102").Value = "Segunda"
$ws.Range("J102").Value = 300
$ws.Range("K102").Value = 1500
$ws.Range("L102").Value = 1500
$ws.Range("M102").Value = 1500
$ws.Range("N102").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O102").Value = "Región de Ñuble"
$ws.Range("P102").Value = 1500
$ws.Range("Q102").Value = 1
$ws.Range("R102").Value = "Hortaliza"
